$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 corresponds to "Ano" = 2025 (faturamento anual update)
$ws.Range("B9").Value = 4166784.62
$ws.Range("C9").Value = 651076.24
$ws.Range("D9").Value = 4817860.86
$ws.Range("E9").Value = 13.51380330232285
$ws.Range("F9").Value = 86.48619669767714
$ws.Range("G9").Value = -37.0767278197978
$ws.Range("H9").Value = -24.75362925756699
$ws.Range("I9").Value = 41910
$ws.Range("J9").Value = 1799
$ws.Range("K9").Value = 43709
$ws.Range("L9").Value = 30308
$ws.Range("M9").Value = 158.9633383925036
$ws.Range("N9").Value = 8.527508842824695

$wb.Save()
